$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells keep their original text representation (inline numeric/percent strings)
# by forcing Text number format before assigning the new values.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "308.42"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "2.12%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "36.21"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "3.17%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.104"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.14%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08125"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "2.31%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.945"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.39%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "4.186"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "3.91%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "7.785"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.81%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9306"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "0.69%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1391"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "15.12%"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "4.87%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09262"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-1.67%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03387"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-4.36%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09860"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.17%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001419"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "2.48%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005750"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-1.79%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.622"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "3.60%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "1.10%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3438"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.12%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1349"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "4.57%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.899"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-2.84%"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "1.42%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04511"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.32%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001218"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.37%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004870"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "6.61%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001241"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.78%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02005"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "5.74%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04954"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "5.07%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007650"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "0.77%"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "7.53%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1385"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "4.64%"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-0.46%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01134"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "1.51%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006424"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "1.61%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000750"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.03%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.001191"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-8.69%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002101"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.03%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002001"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.03%"
